# Update "Full results" and "For plotting" sheets with completed (filtered-age)
# plot results, per commit message "completed plots filtered age".

$wb = $excel.ActiveWorkbook

$wsFull = $wb.Worksheets.Item("Full results")
$wsPlot = $wb.Worksheets.Item("For plotting")

# --- "Full results" sheet (row 2 = NULL MODEL, row 3 = CONDITIONAL MODEL, row 4 = COMPLETE MODEL) ---

# Row 2 - NULL MODEL
$wsFull.Range("C2").Value = 0.891438339859624
$wsFull.Range("D2").Value = 0.108676415936547
$wsFull.Range("E2").Value = 1.00011475579617
$wsFull.Range("J2").Value = 0.108663946118895
$wsFull.Range("K2").Value = 0.102639931948813
$wsFull.Range("L2").Value = 0.010784134150557
$wsFull.Range("M2").Value = 0.0120543342322221
$wsFull.Range("N2").Value = 0.11342406609937

# Row 3 - CONDITIONAL MODEL
$wsFull.Range("F3").Value = 0.890167994015137
$wsFull.Range("G3").Value = 0.102651710475922

# Row 4 - COMPLETE MODEL
$wsFull.Range("H4").Value = 0.879382622322679
$wsFull.Range("I4").Value = 0.0711199949322839
$wsFull.Range("O4").Value = 0.120718280351117

# --- "For plotting" sheet ---

# Row 2 - Sibcorr
$wsPlot.Range("C2").Value = 0.108663946118895
$wsPlot.Range("D2").Value = 0.0492859953086397
$wsPlot.Range("E2").Value = 0.16804189692915
$wsPlot.Range("F2").Value = 948

# Row 3 - IOLIB
$wsPlot.Range("C3").Value = 0.11342406609937
$wsPlot.Range("D3").Value = 0.0385876487127879
$wsPlot.Range("E3").Value = 0.188260483485951
$wsPlot.Range("F3").Value = 948

# Row 4 - IORAD
$wsPlot.Range("C4").Value = 0.120718280351117
$wsPlot.Range("D4").Value = 0.0429775500358876
$wsPlot.Range("E4").Value = 0.198459010666346
$wsPlot.Range("F4").Value = 948
